$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3307615370658065
$ws.Range("C2").Value = 0.3317961784535173
$ws.Range("B3").Value = 37.573121611037
$ws.Range("C3").Value = 37.46721918401393
$ws.Range("B4").Value = 695.1743822837677
$ws.Range("C4").Value = 694.5733146210603
$ws.Range("B5").Value = 60.79136373390147
$ws.Range("C5").Value = 60.58559882360898
$ws.Range("B6").Value = 19110.20413151697
$ws.Range("C6").Value = 19011.04080735726
$ws.Range("B7").Value = -0.1070061023741838
$ws.Range("C7").Value = 0.02777048335857713
$ws.Range("B8").Value = 1439.023794995214
$ws.Range("C8").Value = 1283.983008216596
$ws.Range("B9").Value = 1814.782667572592
$ws.Range("C9").Value = 1813.463598285369
$ws.Range("B10").Value = -0.1162134587908634
$ws.Range("C10").Value = 0.03013246110525187
$ws.Range("B11").Value = 1560.615391370239
$ws.Range("C11").Value = 1393.625188154199
$ws.Range("B12").Value = -3.938413362481739
$ws.Range("C12").Value = -3.937916734372739
$ws.Range("B13").Value = -1.969696105011233
$ws.Range("C13").Value = -1.468862861956752
$ws.Range("B14").Value = -1.91631621329432
$ws.Range("C14").Value = -1.415118491830098
$ws.Range("B15").Value = 1.473401918508976
$ws.Range("C15").Value = 1.999710908679064
